# Reorder the worksheet tabs so "review_info" comes before "hotel_info",
# then add the new "State" column (inserted right after "Hotel_Name", before
# "City") to the "hotel_info" sheet and populate it with "Louisiana".

$wb = $excel.ActiveWorkbook

$reviewInfo = $wb.Worksheets.Item("review_info")
$hotelInfo  = $wb.Worksheets.Item("hotel_info")

# Move review_info so it becomes the first tab (i.e. immediately before
# hotel_info), matching the target tab order.
$reviewInfo.Move($hotelInfo)

# Re-fetch the worksheet handle after the move so subsequent operations act
# on current/live state.
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Insert a new column C (State), shifting City/Zip/... one column to the
# right, and populate header + value.
$hotelInfo.Range("C1").EntireColumn.Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"
